$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time_Interval")
$ws.Activate()
$ws.Rows.Item(2).Insert()
$rng = $ws.Range("A2:D2")
$rng.Font.Name = "Arial Unicode MS"
$rng.Font.Size = 10
$rng.Font.Bold = $false
$rng.Font.Color = 0
$rng.HorizontalAlignment = 1
$rng.VerticalAlignment = -4108
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 100
$ws.Cells.Item(2,4).Value = 6
$ws.Range("F4").Select()
